# Daily attendance processing - rotate "Recorded By" (column G) name lists
# For every row whose column G contains a comma-separated list of recorder
# names, move the first name to the end of the list (left rotation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

# Column G is the "Recorded By" column.
$col = 7

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $text = [string]$val
    if ($text -notmatch ",") { continue }

    $parts = $text -split ",\s*"
    if ($parts.Count -lt 2) { continue }

    $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
    $cell.Value = $rotated
}
